# Apply the two changes recorded in the commit:
#  1. slide5's table switches to table style {A0EB0313-9446-479B-8B5D-CC8B84DDBD56}
#  2. the deck's theme colour scheme (currently the "Integral"/"Red Violet"
#     palette living in ppt/theme/theme2.xml, the theme actually bound to
#     the slide master / presentation) is swapped back to the stock
#     "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{A0EB0313-9446-479B-8B5D-CC8B84DDBD56}")

# --- 2. Theme colours -------------------------------------------------------
# COM ColorScheme items are addressed 1..12 in this fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# .RGB uses the usual OLE COLORREF packing (R + G*256 + B*65536), i.e. a
# "BGR" integer, so hex RRGGBB maps to R + G*256 + B*65536.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}
